# Tidy and polish: add a new "Physics Nut" error code row to the Error Codes sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 27

$ws.Cells.Item($row, 2).Value = 25
$ws.Cells.Item($row, 3).Value = "Physics Nut does not have a collision shape when ready called"
$ws.Cells.Item($row, 4).Value = "PhysicsNut"
$ws.Cells.Item($row, 5).Value = "Ready"
